$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: HFN404CT-ND / CONN FPC/FFC 4POS .5MM SMD GOLD / 4 pin touch connector ---
# (values are written in the same order the original authoring tool produced
# the shared-string table: A6, then C6, then B6)
$ws.Range("A6").Value = "HFN404CT-ND"
$ws.Range("C6").Value = "4 pin touch connector"
$ws.Range("B6").Value = "CONN FPC/FFC 4POS .5MM SMD GOLD"

# --- Row 7: 296-27010-1-ND / DC-DC conv ---
$ws.Range("A7").Value = "296-27010-1-ND"
$ws.Range("B7").Value = "DC-DC conv"

# --- Row 8: BSS138CT-ND / MOSFET N-CH 50V 220MA SOT-23 ---
$ws.Range("A8").Value = "BSS138CT-ND"
$ws.Range("B8").Value = "MOSFET N-CH 50V 220MA SOT-23"

# Re-use the existing "Pnr/part number" style (style index 1, Times 12pt) on the
# cells that need it, by copying the format from an already-styled cell
# instead of re-declaring fonts (which would create a brand-new style entry).
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the row heights used by the rest of the table.
$ws.Rows.Item(6).RowHeight = 14
$ws.Rows.Item(7).RowHeight = 14
$ws.Rows.Item(8).RowHeight = 14

# New column widths.
$ws.Columns.Item(1).ColumnWidth = 17.856026785714285
$ws.Columns.Item(2).ColumnWidth = 18.711495535714285

# Move the active selection below the newly-added rows.
[void]$ws.Range("A9").Select()
